# Applies updated "eficiencia" results across the workbook:
#  - Info: new Objetivo / Tiempo totals
#  - Activados: new set of (Proceso, Tiempo) samples, now 19 rows (every 20 min up to 360)
#  - Operando: Proceso id changed from 4 to 1 for every existing (Tiempo) sample
#  - Contaminantes: new Z / Concentracion values per contaminant

$wb = $excel.ActiveWorkbook

# --- Sheet "Info" ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 1911777828336.237
$wsInfo.Range("B2").Value = 2.174000024795532

# --- Sheet "Activados" ---
$wsAct = $wb.Worksheets.Item("Activados")
$actRows = 19
$arrAct = New-Object 'object[,]' $actRows,2
for ($i = 0; $i -lt $actRows; $i++) {
    $arrAct[$i,0] = 1
    $arrAct[$i,1] = $i * 20
}
$wsAct.Range("A2:B20").Value = $arrAct

# --- Sheet "Operando" ---
$wsOp = $wb.Worksheets.Item("Operando")
$opRows = 365
$arrOp = New-Object 'object[,]' $opRows,2
for ($i = 0; $i -lt $opRows; $i++) {
    $arrOp[$i,0] = 1
    $arrOp[$i,1] = $i
}
$wsOp.Range("A2:B366").Value = $arrOp

# --- Sheet "Contaminantes" ---
$wsCont = $wb.Worksheets.Item("Contaminantes")
$arrCont = New-Object 'object[,]' 5,2
$arrCont[0,0] = 1549768444560.001
$arrCont[0,1] = 57.47700000000005
$arrCont[1,0] = 80215758000
$arrCont[1,1] = 2.975
$arrCont[2,0] = 105655612680
$arrCont[2,1] = 3.918500000000001
$arrCont[3,0] = 386496.2392776002
$arrCont[3,1] = 0.00001433417000000001
$arrCont[4,0] = 176137626600
$arrCont[4,1] = 6.532500000000002
$wsCont.Range("B2:C6").Value = $arrCont
